$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 50,4
$data[0,0] = 0.1685543805360794
$data[0,1] = 0.9514577984809875
$data[0,2] = 0.01198538299649954
$data[0,3] = 0.9987359642982483
$data[1,0] = 0.05668577179312706
$data[1,1] = 0.9885275363922119
$data[1,2] = 0.003683493472635746
$data[1,3] = 0.9991783499717712
$data[2,0] = 0.03427432104945183
$data[2,1] = 0.9928836226463318
$data[2,2] = 0.001357663189992309
$data[2,3] = 0.9994943737983704
$data[3,0] = 0.02565275318920612
$data[3,1] = 0.9932286739349365
$data[3,2] = 0.001237227465026081
$data[3,3] = 0.9994943737983704
$data[4,0] = 0.01652182638645172
$data[4,1] = 0.9966143369674683
$data[4,2] = 0.001224466715939343
$data[4,3] = 0.9994943737983704
$data[5,0] = 0.01606627926230431
$data[5,1] = 0.9965065121650696
$data[5,2] = 0.00120731454808265
$data[5,3] = 0.9994943737983704
$data[6,0] = 0.01491070538759232
$data[6,1] = 0.9967653155326843
$data[6,2] = 0.0006554255960509181
$data[6,3] = 0.9994943737983704
$data[7,0] = 0.01238577254116535
$data[7,1] = 0.997218132019043
$data[7,2] = 0.0003809008921962231
$data[7,3] = 0.9998735785484314
$data[8,0] = 0.01169418264180422
$data[8,1] = 0.9971534609794617
$data[8,2] = 0.0006371518829837441
$data[8,3] = 0.9994943737983704
$data[9,0] = 0.009791713207960129
$data[9,1] = 0.9978435039520264
$data[9,2] = 0.0005319733172655106
$data[9,3] = 0.9994943737983704
$data[10,0] = 0.008975117467343807
$data[10,1] = 0.9980807304382324
$data[10,2] = 0.0006187378894537687
$data[10,3] = 0.9994943737983704
$data[11,0] = 0.009000315330922604
$data[11,1] = 0.9980160593986511
$data[11,2] = 0.0006525839562527835
$data[11,3] = 0.9994943737983704
$data[12,0] = 0.008788539096713066
$data[12,1] = 0.9980160593986511
$data[12,2] = 0.0004736995324492455
$data[12,3] = 0.9994943737983704
$data[13,0] = 0.007557388860732317
$data[13,1] = 0.9984904527664185
$data[13,2] = 0.0004759454459417611
$data[13,3] = 0.9994943737983704
$data[14,0] = 0.007149933371692896
$data[14,1] = 0.9985336065292358
$data[14,2] = 0.0003978715976700187
$data[14,3] = 0.9994943737983704
$data[15,0] = 0.008431801572442055
$data[15,1] = 0.9982532858848572
$data[15,2] = 0.0004011470009572804
$data[15,3] = 0.9994943737983704
$data[16,0] = 0.007538223639130592
$data[16,1] = 0.9984473586082458
$data[16,2] = 0.0003688928845804185
$data[16,3] = 0.9998735785484314
$data[17,0] = 0.007409240584820509
$data[17,1] = 0.9987276792526245
$data[17,2] = 0.0003067189536523074
$data[17,3] = 0.9999368190765381
$data[18,0] = 0.007439374923706055
$data[18,1] = 0.9985767006874084
$data[18,2] = 0.0004118172801099718
$data[18,3] = 0.999620795249939
$data[19,0] = 0.007145338226109743
$data[19,1] = 0.9985982775688171
$data[19,2] = 0.0003451149677857757
$data[19,3] = 0.9999368190765381
$data[20,0] = 0.008178493939340115
$data[20,1] = 0.9984688758850098
$data[20,2] = 0.0002993031812366098
$data[20,3] = 0.9999368190765381
$data[21,0] = 0.009255361743271351
$data[21,1] = 0.9982532858848572
$data[21,2] = 0.0002953959046863019
$data[21,3] = 0.9999368190765381
$data[22,0] = 0.007924159988760948
$data[22,1] = 0.9984257817268372
$data[22,2] = 0.0002573285310063511
$data[22,3] = 0.9999368190765381
$data[23,0] = 0.006933571305125952
$data[23,1] = 0.9985982775688171
$data[23,2] = 0.0002971312787849456
$data[23,3] = 0.9999368190765381
$data[24,0] = 0.007368875201791525
$data[24,1] = 0.9984904527664185
$data[24,2] = 0.0001980973029276356
$data[24,3] = 1
$data[25,0] = 0.007783030159771442
$data[25,1] = 0.9985336065292358
$data[25,2] = 0.0002116225368808955
$data[25,3] = 1
$data[26,0] = 0.007282840553671122
$data[26,1] = 0.9986198544502258
$data[26,2] = 0.0003015203110408038
$data[26,3] = 0.9998735785484314
$data[27,0] = 0.007369357626885176
$data[27,1] = 0.9985120296478271
$data[27,2] = 0.0002253529528388754
$data[27,3] = 0.9999368190765381
$data[28,0] = 0.0076838294044137
$data[28,1] = 0.9985767006874084
$data[28,2] = 0.0003962049377150834
$data[28,3] = 0.9999368190765381
$data[29,0] = 0.006201847456395626
$data[29,1] = 0.9989001750946045
$data[29,2] = 0.0001463020307710394
$data[29,3] = 0.9999368190765381
$data[30,0] = 0.007554376497864723
$data[30,1] = 0.9984688758850098
$data[30,2] = 0.0002316518366569653
$data[30,3] = 0.9999368190765381
$data[31,0] = 0.008164077065885067
$data[31,1] = 0.9984473586082458
$data[31,2] = 0.0002044559660134837
$data[31,3] = 0.9999368190765381
$data[32,0] = 0.007789433002471924
$data[32,1] = 0.9984473586082458
$data[32,2] = 0.0001904767705127597
$data[32,3] = 0.9999368190765381
$data[33,0] = 0.006942082196474075
$data[33,1] = 0.9985982775688171
$data[33,2] = 0.0001555290655232966
$data[33,3] = 1
$data[34,0] = 0.006255006417632103
$data[34,1] = 0.9987923502922058
$data[34,2] = 0.0001506775297457352
$data[34,3] = 1
$data[35,0] = 0.006362416781485081
$data[35,1] = 0.9986845254898071
$data[35,2] = 0.000153501721797511
$data[35,3] = 0.9999368190765381
$data[36,0] = 0.007286135107278824
$data[36,1] = 0.9985982775688171
$data[36,2] = 0.0002799557114485651
$data[36,3] = 0.9999368190765381
$data[37,0] = 0.007281627040356398
$data[37,1] = 0.9986198544502258
$data[37,2] = 0.0002107537002302706
$data[37,3] = 0.9999368190765381
$data[38,0] = 0.007905551232397556
$data[38,1] = 0.9984688758850098
$data[38,2] = 0.0001422477071173489
$data[38,3] = 1
$data[39,0] = 0.006917495746165514
$data[39,1] = 0.9985767006874084
$data[39,2] = 0.0001651209458941594
$data[39,3] = 1
$data[40,0] = 0.00666859420016408
$data[40,1] = 0.9986414313316345
$data[40,2] = 0.000148038103361614
$data[40,3] = 1
$data[41,0] = 0.006477537099272013
$data[41,1] = 0.9986414313316345
$data[41,2] = 0.0001918048219522461
$data[41,3] = 1
$data[42,0] = 0.009205807000398636
$data[42,1] = 0.9983395338058472
$data[42,2] = 0.0001872269494924694
$data[42,3] = 0.9999368190765381
$data[43,0] = 0.007009216118603945
$data[43,1] = 0.9986414313316345
$data[43,2] = 0.0001148005758295767
$data[43,3] = 1
$data[44,0] = 0.007330498192459345
$data[44,1] = 0.9985120296478271
$data[44,2] = 0.0001330219383817166
$data[44,3] = 0.9999368190765381
$data[45,0] = 0.006888146977871656
$data[45,1] = 0.9985551834106445
$data[45,2] = 0.0001357696455670521
$data[45,3] = 0.9999368190765381
$data[46,0] = 0.007587114814668894
$data[46,1] = 0.9984904527664185
$data[46,2] = 0.0001083463721442968
$data[46,3] = 1
$data[47,0] = 0.006127714645117521
$data[47,1] = 0.9987276792526245
$data[47,2] = 0.0001025825113174506
$data[47,3] = 1
$data[48,0] = 0.006919529289007187
$data[48,1] = 0.9985982775688171
$data[48,2] = 0.00009808960749069229
$data[48,3] = 1
$data[49,0] = 0.007372591178864241
$data[49,1] = 0.9984688758850098
$data[49,2] = 0.0001248484913958237
$data[49,3] = 1

$ws.Range("A2:D51").Value = $data

